# "Creacion de consola de planificacion"
# Populate the "Reporte atrasos" (sheet2) and "Reporte planificacion" (sheet3)
# worksheets with the new planning-console rows, then reset view/selection
# state so "Sheet" (the first tab) is the active sheet with A1 selected
# everywhere.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 2: "Reporte atrasos" - add a new data row (row 4)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A4").Value = 1

$ws2.Range("B4:C4").NumberFormat = "yyyy\-mm\-dd"
$ws2.Range("B4").Value = "2017-03-30"
$ws2.Range("C4").Value = "2017-04-10"

# widen column C to match column B (both 30 chars wide)
$ws2.Columns.Item(3).ColumnWidth = 29.2

# ---------------------------------------------------------------------------
# Sheet 3: "Reporte planificación" - rewrite row 4 and add rows 5-7
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# apply the date display format to every date column BEFORE writing values,
# so Excel doesn't have to invent its own default date style on entry
$ws3.Range("B4:C7").NumberFormat = "yyyy\-mm\-dd"
$ws3.Range("E4:F7").NumberFormat = "yyyy\-mm\-dd"
$ws3.Range("H4:I7").NumberFormat = "yyyy\-mm\-dd"
$ws3.Range("K4:L7").NumberFormat = "yyyy\-mm\-dd"

# row 5 (D5/M5) carried leftover blank-cell styling from the template that
# must be cleared before the new values are written
$ws3.Range("D5").ClearFormats()
$ws3.Range("M5").ClearFormats()

# --- row 4 (overwrite existing demo values) ---
$ws3.Range("A4").Value = 1
$ws3.Range("B4").Value = "2017-03-23"
$ws3.Range("C4").Value = "2017-03-28"
$ws3.Range("D4").Value = 9
$ws3.Range("E4").Value = "2017-03-29"
$ws3.Range("F4").Value = "2017-04-03"
$ws3.Range("G4").Value = 6
$ws3.Range("H4").Value = "2017-04-04"
$ws3.Range("I4").Value = "2017-04-07"
$ws3.Range("J4").Value = 7
$ws3.Range("K4").Value = "2017-04-10"
$ws3.Range("L4").Value = "2017-04-10"
$ws3.Range("M4").Value = "8 ;12 ;13"

# --- row 5 ---
$ws3.Range("A5").Value = 2
$ws3.Range("B5").Value = "2017-03-23"
$ws3.Range("C5").Value = "2017-03-29"
$ws3.Range("D5").Value = 5
$ws3.Range("E5").Value = "2017-03-30"
$ws3.Range("F5").Value = "2017-04-05"
$ws3.Range("G5").Value = 10
$ws3.Range("H5").Value = "2017-04-06"
$ws3.Range("I5").Value = "2017-04-12"
$ws3.Range("J5").Value = 11
$ws3.Range("K5").Value = "2017-04-13"
$ws3.Range("L5").Value = "2017-04-20"

$ws3.Range("M5").NumberFormat = "@"
$ws3.Range("M5").Value = "13"
$ws3.Range("M5").ClearFormats()

# --- row 6 (brand new row) ---
$ws3.Range("A6").Value = 3
$ws3.Range("B6").Value = "2017-03-23"
$ws3.Range("C6").Value = "2017-04-07"
$ws3.Range("D6").Value = 1
$ws3.Range("E6").Value = "2017-04-10"
$ws3.Range("F6").Value = "2017-04-27"
$ws3.Range("G6").Value = 10
$ws3.Range("H6").Value = "2017-04-28"
$ws3.Range("I6").Value = "2017-05-16"
$ws3.Range("J6").Value = 11
$ws3.Range("K6").Value = "2017-05-17"
$ws3.Range("L6").Value = "2017-05-29"

$ws3.Range("M6").NumberFormat = "@"
$ws3.Range("M6").Value = "13"
$ws3.Range("M6").ClearFormats()

# --- row 7 (brand new row) ---
$ws3.Range("A7").Value = 4
$ws3.Range("B7").Value = "2017-03-29"
$ws3.Range("C7").Value = "2017-04-05"
$ws3.Range("D7").Value = 9
$ws3.Range("E7").Value = "2017-04-06"
$ws3.Range("F7").Value = "2017-04-13"
$ws3.Range("G7").Value = 6
$ws3.Range("H7").Value = "2017-04-17"
$ws3.Range("I7").Value = "2017-04-25"
$ws3.Range("J7").Value = 11
$ws3.Range("K7").Value = "2017-04-26"
$ws3.Range("L7").Value = "2017-05-02"

$ws3.Range("M7").NumberFormat = "@"
$ws3.Range("M7").Value = "13"
$ws3.Range("M7").ClearFormats()

# widen every data column (2-13) to match the console layout
for ($c = 2; $c -le 13; $c++) {
    $ws3.Columns.Item($c).ColumnWidth = 29.2
}

# ---------------------------------------------------------------------------
# View state: console opens on the first sheet, with A1 selected everywhere
# ---------------------------------------------------------------------------
[void]$ws3.Range("A1").Select()

$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("A1").Select()

[void]$ws2.Range("A1").Select()

[void]$ws1.Activate()
[void]$ws1.Range("A1").Select()
